$d = $word.ActiveDocument

# 1. Replace the title text "FÉRIAS" with the new heading text.
$d.Content.Find.Execute("FÉRIAS", $true, $false, $false, $false, $false, $true, 1, $false, `
    "CHECK LIST – GUIA PARA ENTREVISTA", 2)

$para = $d.Paragraphs(1)

# 2. Paragraph formatting: a left tab stop at 1800 twips (90 pt) and centered text.
$para.Range.ParagraphFormat.TabStops.Add(90)
$para.Range.ParagraphFormat.Alignment = 1

# 3. Character formatting for the whole paragraph (run + paragraph mark): bold, 16pt (32 half-points),
#    including the complex-script mirrors so bCs/szCs are emitted too.
$para.Range.Font.Bold = $true
$para.Range.Font.BoldBi = $true
$para.Range.Font.Size = 16
$para.Range.Font.SizeBi = 16

# 4. Add the "_GoBack" bookmark collapsed at the very start of the document/paragraph.
#    Creating a collapsed bookmark exactly at offset 0 directly tends to have it swallow the
#    whole paragraph's run on save, so first anchor it one character in (after a throwaway
#    placeholder), then delete the placeholder - the bookmark stays collapsed at the true start.
$lead = $d.Range(0, 0)
$lead.InsertBefore("X")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range(0, 1)
$placeholder.Delete()
